$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "nemad" (symbol) column L for all data rows (2-45) from the
# company name "نفت سپاهان" to the actual ticker symbol "شسپا".
for ($r = 2; $r -le 45; $r++) {
    $ws.Cells.Item($r, 12).Value = "شسپا"
}

# Column L width was auto-fit by Excel after the text changed (renders as
# width 10 in the saved XML).
$ws.Columns.Item(12).ColumnWidth = 9.14

# Record the active selection left on L5 before the file was saved.
$ws.Range("L5").Select()
